{"js": "// Update the date line and each multiplication expression in the practice\n// sheet table. Every \"before\" string is unique in the document, so a plain\n// search+replace per pair is unambiguous.\nconst replacements = [\n  [\"2024-03-22 Friday\", \"2024-03-23 Saturday\"],\n  [\"255\u00d74=\", \"160\u00d76=\"],\n  [\"267\u00d74=\", \"342\u00d79=\"],\n  [\"217\u00d76=\", \"415\u00d78=\"],\n  [\"954\u00d79=\", \"856\u00d74=\"],\n  [\"529\u00d76=\", \"731\u00d74=\"],\n  [\"261\u00d75=\", \"424\u00d72=\"],\n  [\"672\u00d73=\", \"111\u00d74=\"],\n  [\"113\u00d73=\", \"863\u00d75=\"],\n  [\"318\u00d78=\", \"110\u00d76=\"],\n  [\"222\u00d76=\", \"992\u00d79=\"],\n  [\"345\u00d74=\", \"245\u00d72=\"],\n  [\"101\u00d75=\", \"447\u00d75=\"],\n  [\"198\u00d75=\", \"980\u00d73=\"],\n  [\"846\u00d77=\", \"750\u00d77=\"],\n  [\"723\u00d78=\", \"726\u00d78=\"],\n  [\"580\u00d76=\", \"349\u00d79=\"],\n  [\"484\u00d73=\", \"118\u00d79=\"],\n  [\"394\u00d78=\", \"579\u00d72=\"],\n  [\"177\u00d79=\", \"481\u00d79=\"],\n  [\"609\u00d75=\", \"875\u00d78=\"],\n  [\"202\u00d74=\", \"983\u00d72=\"],\n  [\"206\u00d75=\", \"269\u00d79=\"],\n  [\"212\u00d76=\", \"423\u00d72=\"],\n  [\"231\u00d77=\", \"509\u00d74=\"],\n  [\"204\u00d72=\", \"437\u00d74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for '\" + before + \"' but found \" + results.items.length\n    );\n  }\n\n  results.items[0].insertText(after, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the date line and each multiplication expression in the practice\n# sheet table. Every \"before\" string is unique in the document, so a plain\n# Find/Replace per pair (ReplaceOne) is unambiguous and safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-03-22 Friday', '2024-03-23 Saturday'),\n    @('255\u00d74=', '160\u00d76='),\n    @('267\u00d74=', '342\u00d79='),\n    @('217\u00d76=', '415\u00d78='),\n    @('954\u00d79=', '856\u00d74='),\n    @('529\u00d76=', '731\u00d74='),\n    @('261\u00d75=', '424\u00d72='),\n    @('672\u00d73=', '111\u00d74='),\n    @('113\u00d73=', '863\u00d75='),\n    @('318\u00d78=', '110\u00d76='),\n    @('222\u00d76=', '992\u00d79='),\n    @('345\u00d74=', '245\u00d72='),\n    @('101\u00d75=', '447\u00d75='),\n    @('198\u00d75=', '980\u00d73='),\n    @('846\u00d77=', '750\u00d77='),\n    @('723\u00d78=', '726\u00d78='),\n    @('580\u00d76=', '349\u00d79='),\n    @('484\u00d73=', '118\u00d79='),\n    @('394\u00d78=', '579\u00d72='),\n    @('177\u00d79=', '481\u00d79='),\n    @('609\u00d75=', '875\u00d78='),\n    @('202\u00d74=', '983\u00d72='),\n    @('206\u00d75=', '269\u00d79='),\n    @('212\u00d76=', '423\u00d72='),\n    @('231\u00d77=', '509\u00d74='),\n    @('204\u00d72=', '437\u00d74=')\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute(\n        $findText,   # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $replaceText,# ReplaceWith\n        2            # Replace (wdReplaceOne)\n    )\n\n    if (-not $found) {\n        throw \"Find.Execute did not find expected text: $findText\"\n    }\n}\n"}
